$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp on row 4 (tiny precision correction from the data refresh)
$ws.Range("A4").Value = 45804.43702679398

# Append new row 5 with the latest price entry
$ws.Range("A5").Value = 45804.44030289967
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("B5").Value = "EVOWHEY PROTEIN"
$ws.Range("C5").Value = "2Kg"
$ws.Range("D5").Value = "37,90€"
